$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 620.9
$ws.Range("I2").Value = 330.2
$ws.Range("J2").Value = 911.6
$ws.Range("K2").Value = 330.2
$ws.Range("L2").Value = 911.6
$ws.Range("M2").Value = -217.2
$ws.Range("N2").Value = -1137.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 75998.5
$ws.Range("J87").Value = 75998.5
$ws.Range("L87").Value = 75998.5
$ws.Range("N87").Value = -78494.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 75998.5
$ws.Range("J90").Value = 75998.5
$ws.Range("L90").Value = 227995.5
$ws.Range("N90").Value = -240475.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1421.3846
$ws.Range("I98").Value = 1435.6666
$ws.Range("K98").Value = 1435.6666
$ws.Range("M98").Value = 62.33339999999998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 23675.176
$ws.Range("I106").Value = 23675.176
$ws.Range("K106").Value = 23675.176
$ws.Range("M106").Value = -23044.176

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 2454.2942
$ws.Range("I107").Value = 2541.2666
$ws.Range("K107").Value = 2541.2666
$ws.Range("M107").Value = -621.2665999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1842.8572
$ws.Range("J112").Value = 2000
$ws.Range("L112").Value = 6000
$ws.Range("N112").Value = -8216

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1421.3846
$ws.Range("I122").Value = 1435.6666
$ws.Range("K122").Value = 4306.9998
$ws.Range("M122").Value = -1856.9998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1701.0667
$ws.Range("I132").Value = 1701.0667
$ws.Range("K132").Value = 5103.2001
$ws.Range("M132").Value = -2573.2001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2045.3455
$ws.Range("I138").Value = 1733.6154
$ws.Range("K138").Value = 5200.8462
$ws.Range("M138").Value = -60.84619999999995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3200.875
$ws.Range("I45").Value = 2003.5
$ws.Range("K45").Value = 2003.5
$ws.Range("M45").Value = -1626.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 646.61536
$ws.Range("I97").Value = 668.1667
$ws.Range("J97").Value = 388
$ws.Range("K97").Value = 668.1667
$ws.Range("L97").Value = 388
$ws.Range("M97").Value = -172.1667
$ws.Range("N97").Value = -1380

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1579.8948
$ws.Range("I94").Value = 901.4286
$ws.Range("J94").Value = 3479.6
$ws.Range("K94").Value = 901.4286
$ws.Range("L94").Value = 3479.6
$ws.Range("M94").Value = -450.4286
$ws.Range("N94").Value = -4381.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3984
$ws.Range("J134").Value = 3984
$ws.Range("L134").Value = 11952
$ws.Range("N134").Value = -17022

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 5500
$ws.Range("I4").Value = 5500
$ws.Range("K4").Value = 5500
$ws.Range("M4").Value = -5388

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 5000
$ws.Range("I39").Value = 5000
$ws.Range("K39").Value = 5000
$ws.Range("M39").Value = -4609

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H49").Value = 5000
$ws.Range("I49").Value = 5000
$ws.Range("K49").Value = 5000
$ws.Range("M49").Value = -4818

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 12424.792
$ws.Range("I99").Value = 7709.154
$ws.Range("K99").Value = 7709.154
$ws.Range("M99").Value = -6211.154

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 12424.792
$ws.Range("I126").Value = 7709.154
$ws.Range("K126").Value = 23127.462
$ws.Range("M126").Value = -20657.462

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2973.5
$ws.Range("I132").Value = 2631.3333
$ws.Range("K132").Value = 7893.999899999999
$ws.Range("M132").Value = -5363.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 59333.117
$ws.Range("I107").Value = 388.25
$ws.Range("J107").Value = 77470
$ws.Range("K107").Value = 1164.75
$ws.Range("L107").Value = 232410
$ws.Range("M107").Value = 755.25
$ws.Range("N107").Value = -236250

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1783
$ws.Range("I131").Value = 1605
$ws.Range("K131").Value = 4815
$ws.Range("M131").Value = 225

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 145
$ws.Range("I19").Value = 140
$ws.Range("K19").Value = 140
$ws.Range("M19").Value = 148

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 14142.857
$ws.Range("J20").Value = 14142.857
$ws.Range("L20").Value = 14142.857
$ws.Range("N20").Value = -14632.857

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 16444.445
$ws.Range("J24").Value = 16444.445
$ws.Range("L24").Value = 16444.445
$ws.Range("N24").Value = -16790.445

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 3760000
$ws.Range("I35").Value = 3346666.8
$ws.Range("J35").Value = 5000000
$ws.Range("K35").Value = 3346666.8
$ws.Range("L35").Value = 5000000
$ws.Range("M35").Value = -3346368.8
$ws.Range("N35").Value = -5000596

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3450.2104
$ws.Range("I80").Value = 2991.6667
$ws.Range("J80").Value = 3661.8462
$ws.Range("K80").Value = 2991.6667
$ws.Range("L80").Value = 3661.8462
$ws.Range("M80").Value = -1993.6667
$ws.Range("N80").Value = -5657.8462

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3450.2104
$ws.Range("I83").Value = 2991.6667
$ws.Range("J83").Value = 3661.8462
$ws.Range("K83").Value = 14958.3335
$ws.Range("L83").Value = 18309.231
$ws.Range("M83").Value = -9966.333500000001
$ws.Range("N83").Value = -28293.231

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2906.8125
$ws.Range("I97").Value = 2374.875
$ws.Range("J97").Value = 3438.75
$ws.Range("K97").Value = 2374.875
$ws.Range("L97").Value = 3438.75
$ws.Range("M97").Value = -1878.875
$ws.Range("N97").Value = -4430.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1714.1428
$ws.Range("I132").Value = 1499.8334
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 4499.5002
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -1969.5002
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 939.1667
$ws.Range("I55").Value = 841.75
$ws.Range("J55").Value = 1134
$ws.Range("K55").Value = 841.75
$ws.Range("L55").Value = 1134
$ws.Range("M55").Value = -668.75
$ws.Range("N55").Value = -1480

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1292.8
$ws.Range("I93").Value = 1116.625
$ws.Range("J93").Value = 1997.5
$ws.Range("K93").Value = 1116.625
$ws.Range("L93").Value = 1997.5
$ws.Range("M93").Value = 131.375
$ws.Range("N93").Value = -4493.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5532
$ws.Range("I132").Value = 5532
$ws.Range("K132").Value = 16596
$ws.Range("M132").Value = -14066

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 30000
$ws.Range("J31").Value = 30000
$ws.Range("L31").Value = 30000
$ws.Range("N31").Value = -30696

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2273
$ws.Range("I122").Value = 3339.8
$ws.Range("J122").Value = 495
$ws.Range("K122").Value = 10019.4
$ws.Range("L122").Value = 1485
$ws.Range("M122").Value = -7569.400000000001
$ws.Range("N122").Value = -6385

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 80000
$ws.Range("J135").Value = 80000
$ws.Range("L135").Value = 80000
$ws.Range("N135").Value = -90140
